# Apply updated case_id / prompt_tokens / completion_tokens values
# to the gemini_outputs worksheet (Sheet1), rows 2-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(A case_id, D prompt_tokens, E completion_tokens)
$data = @{
    2  = @("azd-2_22-cv-02126",    113728, 924)
    3  = @("cand-3_18-cv-04865",   29090,  748)
    4  = @("cand-4_22-cv-02672",   47793,  560)
    5  = @("cand_22_cv_02094",     79090,  579)
    6  = @("cand_23_cv_02560",     69684,  927)
    7  = @("cand_23_cv_03518",     33551,  759)
    8  = @("cand_24_cv_03170",     26991,  626)
    9  = @("cand_24_cv_04196",     12688,  568)
    10 = @("cand_3_22-cv-00956",   21697,  549)
    11 = @("casd_3_23-cv-01216",   34901,  759)
    12 = @("ctd-3-23-cv-01035",    66233,  683)
    13 = @("dcd-1_23-cv-02055",    38767,  591)
    14 = @("dde_ 23_cv_1466",      35588,  812)
    15 = @("dde_21_cv_55",         47422,  661)
    16 = @("flsd-1_23-cv-23139",   16873,  1047)
    17 = @("ilnd-1-21-cv-04349",   32712,  722)
    18 = @("mad-1-21-cv-10933",    17417,  628)
    19 = @("mied-4-23-cv-13132",   67564,  894)
    20 = @("nysd_20_cv_04494",     54158,  1376)
    21 = @("nysd_22-cv-07111",     28976,  902)
    22 = @("nysd_22_cv_10292",     24881,  785)
    23 = @("nysd_23_cv_9476",      16625,  335)
    24 = @("nysd_24_cv_310",       45369,  943)
    25 = @("txnd-4_24-cv-00673",   58520,  537)
    26 = @("txsd-4-21-cv-02473",   70610,  588)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
}
